# Update cryptocurrency price/volume figures (Sheet1 columns D and E)
# to match the latest scrape, preserving each cell as plain text
# (matching the inline-string storage of the source workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.398.27"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.659.70"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.18%  "
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").Value = "2.657.86"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "3.148.13"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").Value = "72.293.30"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "2.653.24"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "0.0₃0954"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "498.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("E42").Value = "  -6.82%  "
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.329"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.552"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").Value = "  +0.35%  "
